$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the old row 8 ("Starlight Traveler"),
# shifting all subsequent rows down by two.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# New song 1 -> row 8
$ws.Range("A8").Value = "回不去的夏天"
$ws.Range("B8").Value = "夏日入侵企画"
$ws.Range("C8").Value = "mrY8qdNHcRs"

# New song 2 -> row 9
$ws.Range("A9").Value = "我們都擁有海洋 (嗶哩嗶哩2023畢業歌)"
$ws.Range("B9").Value = "吳青峰"
$ws.Range("C9").Value = "rQOIRBrY7h0"

# Update the saved view/selection state to match the edited workbook.
$ws.Range("E10").Select()

Write-Output "done"
